$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column E (header "D") values: reveal/hide different cells ---
# Row 2: was missing -> now has a value
$ws.Range("E2").Value = -7.2

# Row 6: had a value -> now missing
$ws.Range("E6").Value = $null

# Row 12: was missing -> now has a value
$ws.Range("E12").Value = -5.3

# Row 14: had a value -> now missing
$ws.Range("E14").Value = $null

# Row 20: was missing -> now has a value
$ws.Range("E20").Value = -7.2

# Row 21: was missing -> now has a value
$ws.Range("E21").Value = -8.699999999999999

# Row 23: had a value -> now missing
$ws.Range("E23").Value = $null

# Row 24: had a value -> now missing
$ws.Range("E24").Value = $null

# --- Remove two data rows entirely (sheet shrinks from 34 to 32 data rows) ---
# Delete "SC 92" (row 28) first so "RM 232" keeps its row index of 26.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- After the deletions, the "SC ..." rows occupy rows 26-33. Update column C ---
# (header "B") and the remaining column E (header "D") cells to match the new
# reveal/hide pattern of missing values for this re-sampled dataset.
$ws.Range("C26").Value = 10.8        # SC 5:   was missing -> now has a value
$ws.Range("C27").Value = $null       # SC 101: had a value -> now missing
$ws.Range("C28").Value = $null       # SC 105: had a value -> now missing
$ws.Range("C29").Value = 11.2        # SC 119: was missing -> now has a value
$ws.Range("C30").Value = 11.4        # SC 120: was missing -> now has a value
$ws.Range("C31").Value = $null       # SC 132: had a value -> now missing
$ws.Range("E31").Value = -8.1        # SC 132: was missing -> now has a value
$ws.Range("C32").Value = $null       # SC 193: had a value -> now missing
$ws.Range("E33").Value = -10.7       # SC 232: was missing -> now has a value
